$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bumil_TPK - Bumil")

# Shorten the "Desa" column header (T1) from "Desa di Kec. Bulakamba" to "Desa"
$ws.Range("T1").Value = "Desa"

# Match the author's resulting selection (cell T2 was clicked/selected last)
$ws.Range("T2").Select() | Out-Null
